# faturamento_diario.xlsx update:
# Insert 3 new daily rows (dias 11, 12, 13) for Julho/2025 at the top of the
# June block (current row 12), pushing the existing June/May/Abril rows down
# by 3. This grows the sheet from 102 to 105 rows (A1:E105) and brings the
# data for period 07/2025 from 10 days up to 13 days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 12 (start of the June data),
# one at a time so each new row lands directly beneath the prior new one,
# keeping July's data contiguous and in day order (11, 12, 13).
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()

# New row 12 -> Dia 11, Julho/2025
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 25398.22
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 2025
$ws.Range("E12").Value = "07/2025"

# New row 13 -> Dia 12, Julho/2025
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 10045.9
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 2025
$ws.Range("E13").Value = "07/2025"

# New row 14 -> Dia 13, Julho/2025
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 2989.9
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 2025
$ws.Range("E14").Value = "07/2025"
